# Auto-applied cell updates from diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue $ws 'D2' '29.076.24'
Set-TextValue $ws 'E2' '  -0.53%  '
Set-TextValue $ws 'D3' '1.820.61'
Set-TextValue $ws 'E3' '  -0.52%  '
Set-TextValue $ws 'E4' '  +0.46%  '
Set-TextValue $ws 'D5' '233.08'
Set-TextValue $ws 'E5' '  -1.77%  '
Set-TextValue $ws 'D6' '0.5905'
Set-TextValue $ws 'E6' '  -3.11%  '
Set-TextValue $ws 'E7' '  +0.45%  '
Set-TextValue $ws 'D8' '0.2739'
Set-TextValue $ws 'E8' '  -3.22%  '
Set-TextValue $ws 'D9' '0.06789'
Set-TextValue $ws 'E9' '  -4.62%  '
Set-TextValue $ws 'D10' '22.94'
Set-TextValue $ws 'E10' '  -4.42%  '
Set-TextValue $ws 'D11' '0.07510'
Set-TextValue $ws 'E11' '  -1.76%  '
Set-TextValue $ws 'D12' '1.831.95'
Set-TextValue $ws 'E12' '  -0.20%  '
Set-TextValue $ws 'D13' '4.667'
Set-TextValue $ws 'E13' '  -3.19%  '
Set-TextValue $ws 'D14' '0.6233'
Set-TextValue $ws 'E14' '  -2.43%  '
Set-TextValue $ws 'D15' '0.000009345'
Set-TextValue $ws 'E15' '  -6.44%  '
Set-TextValue $ws 'D16' '74.40'
Set-TextValue $ws 'E16' '  -6.76%  '
Set-TextValue $ws 'D17' '28.827.15'
Set-TextValue $ws 'E17' '  -1.30%  '
Set-TextValue $ws 'D18' '5.411'
Set-TextValue $ws 'E18' '  -9.69%  '
Set-TextValue $ws 'D19' '1.005'
Set-TextValue $ws 'E19' '  +0.44%  '
Set-TextValue $ws 'D20' '207.63'
Set-TextValue $ws 'E20' '  -9.81%  '
Set-TextValue $ws 'D21' '11.36'
Set-TextValue $ws 'E21' '  -4.02%  '
Set-TextValue $ws 'D22' '6.769'
Set-TextValue $ws 'E22' '  -4.01%  '
Set-TextValue $ws 'D23' '1.007'
Set-TextValue $ws 'E23' '  +0.39%  '
Set-TextValue $ws 'D24' '153.88'
Set-TextValue $ws 'E24' '  -1.03%  '
Set-TextValue $ws 'D25' '0.1268'
Set-TextValue $ws 'E25' '  -2.36%  '
Set-TextValue $ws 'D26' '7.771'
Set-TextValue $ws 'E26' '  -4.15%  '
Set-TextValue $ws 'D27' '16.24'
Set-TextValue $ws 'E27' '  -3.11%  '
Set-TextValue $ws 'D28' '0.06392'
Set-TextValue $ws 'E28' '  -6.47%  '
Set-TextValue $ws 'D29' '1.408'
Set-TextValue $ws 'E29' '  -4.85%  '
Set-TextValue $ws 'E30' '  -1.92%  '
Set-TextValue $ws 'D31' '3.704'
Set-TextValue $ws 'E31' '  -3.32%  '
Set-TextValue $ws 'D32' '3.670'
Set-TextValue $ws 'E32' '  -4.60%  '
Set-TextValue $ws 'D33' '1.673'
Set-TextValue $ws 'E33' '  -3.69%  '
Set-TextValue $ws 'D34' '1.048'
Set-TextValue $ws 'E34' '  -7.09%  '
Set-TextValue $ws 'D35' '2.532'
Set-TextValue $ws 'E35' '  -0.88%  '
Set-TextValue $ws 'D36' '0.6311'
Set-TextValue $ws 'E36' '  -4.32%  '
Set-TextValue $ws 'D37' '2.755'
Set-TextValue $ws 'E37' '  -0.14%  '
Set-TextValue $ws 'D38' '6.457'
Set-TextValue $ws 'E38' '  -2.29%  '
Set-TextValue $ws 'D39' '0.01702'
Set-TextValue $ws 'E39' '  -3.61%  '
Set-TextValue $ws 'D40' '1.129.99'
Set-TextValue $ws 'E40' '  -8.34%  '
Set-TextValue $ws 'D41' '0.8718'
Set-TextValue $ws 'E41' '  -6.65%  '
Set-TextValue $ws 'D42' '1.006'
Set-TextValue $ws 'E42' '  +0.50%  '
Set-TextValue $ws 'D43' '1.974.15'
Set-TextValue $ws 'E43' '  -0.60%  '
Set-TextValue $ws 'D44' '99.57'
Set-TextValue $ws 'E44' '  -1.33%  '
Set-TextValue $ws 'D45' '60.01'
Set-TextValue $ws 'E45' '  -5.65%  '
Set-TextValue $ws 'D46' '0.00000000113'
Set-TextValue $ws 'E46' '  -3.96%  '
Set-TextValue $ws 'D47' '1.573'
Set-TextValue $ws 'E47' '  -3.76%  '
Set-TextValue $ws 'B48' 'Cronos'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D48' '0.05487'
Set-TextValue $ws 'E48' '  -1.38%  '
Set-TextValue $ws 'B49' 'Mantle'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws 'D49' '0.4524'
Set-TextValue $ws 'E49' '  -0.86%  '
Set-TextValue $ws 'D50' '1.011'
Set-TextValue $ws 'E50' '  +0.80%  '
Set-TextValue $ws 'D51' '8.212'
Set-TextValue $ws 'E51' '  -4.22%  '
